$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the exception-description text in D14 ("5.1. ...")
$ws.Range("D14").Value = "5.1. Informa cliente que já existe um user com o email com o qual se pretendia registar"

# Widen column D to accommodate the longer text
$ws.Columns("D").ColumnWidth = 86.25

# Update the active selection shown when the sheet was last saved
$ws.Range("F17").Select()
